$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update its "through" date label
$ws.Name = "Through 2022-09-10"
$ws.Range("A10").Value = "September (through 09-10)"

# Update September row (row 10) values for each year column B..I
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 17
$ws.Range("D10").Value = 25
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 23
$ws.Range("G10").Value = 32
$ws.Range("H10").Value = 44
$ws.Range("I10").Value = 48

# Update Total row (row 11) values for each year column B..I
$ws.Range("B11").Value = 203
$ws.Range("C11").Value = 398
$ws.Range("D11").Value = 576
$ws.Range("E11").Value = 505
$ws.Range("F11").Value = 372
$ws.Range("G11").Value = 816
$ws.Range("H11").Value = 1114
$ws.Range("I11").Value = 1185
